$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new feature tracker entry as row 42 (row 41 was previously the last row)
$ws.Range("A42").Value = "In-App purchases"
$ws.Range("C42").Value = "1.8.5"
$ws.Range("B42").Value = "I want to be able to have enough money to fund the iOS version of the app"
$ws.Range("D42").Value = "Weston Fiala"
